# Applies the updated Northwind report figures (refreshed data pipeline run)
# across all six report sheets: KPIs, Ventes Mensuelles, Par Catégorie,
# Top Produits, Par Pays and Employés.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# KPIs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("KPIs")
$ws.Range("A2").Value = 7831.599999999999
$ws.Range("E2").Value = 163.1583333333333
$ws.Range("F2").Value = 0.8712871287128713

# ---------------------------------------------------------------------
# Ventes Mensuelles
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ventes Mensuelles")
$ws.Range("C2").Value = 717.1
$ws.Range("E2").Value = 49

$ws.Range("C3").Value = 861.0999999999999
$ws.Range("E3").Value = 52

$ws.Range("C4").Value = 1069.3
$ws.Range("E4").Value = 70

$ws.Range("C5").Value = 2690.3
$ws.Range("E5").Value = 192

$ws.Range("C6").Value = 1418.3
$ws.Range("E6").Value = 103

$ws.Range("C7").Value = 1075.5
$ws.Range("E7").Value = 84

# ---------------------------------------------------------------------
# Par Catégorie
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Par Catégorie")
$ws.Range("B2").Value = 4541.4
$ws.Range("D2").Value = 275

$ws.Range("B3").Value = 3290.2
$ws.Range("C3").Value = 36
$ws.Range("D3").Value = 275

# ---------------------------------------------------------------------
# Top Produits
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Top Produits")
$ws.Range("B2").Value = 4541.4
$ws.Range("C2").Value = 275

$ws.Range("B3").Value = 1847
$ws.Range("C3").Value = 203
$ws.Range("D3").Value = 36

$ws.Range("B4").Value = 1443.2
$ws.Range("C4").Value = 72
$ws.Range("D4").Value = 17

# ---------------------------------------------------------------------
# Par Pays
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Par Pays")
$ws.Range("B2").Value = 7831.6

# ---------------------------------------------------------------------
# Employés (also re-sorted by TotalSales descending -> row order changes)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Employés")

$ws.Range("A2").Value = "Anne Hellung-Larsen"
$ws.Range("B2").Value = 1941.1
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 5

$ws.Range("A3").Value = "Nancy Freehafer"
$ws.Range("B3").Value = 1826.5
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 6

$ws.Range("A4").Value = "Mariya Sergienko"
$ws.Range("B4").Value = 1415.3
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 4

$ws.Range("B5").Value = 950.3

$ws.Range("B6").Value = 658.7

$ws.Range("B7").Value = 457

$ws.Range("B8").Value = 375.5

$ws.Range("B9").Value = 207.2
